$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1976573938506589
$ws.Range("C2").Value = 0.5490483162518301
$ws.Range("J2").Value = 0.01317715959004392
$ws.Range("P2").Value = 0.1376281112737921
$ws.Range("S2").Value = 0.102489019033675
$ws.Range("B3").Value = 0.01015228426395939
$ws.Range("C3").Value = 0.03045685279187817
$ws.Range("J3").Value = 0.03553299492385787
$ws.Range("P3").Value = 0.7131979695431472
$ws.Range("S3").Value = 0.2106598984771574
$ws.Range("J4").Value = 0.09677419354838709
$ws.Range("P4").Value = 0.7634408602150538
$ws.Range("S4").Value = 0.1397849462365591
$ws.Range("B6").Value = 0.07522935779816514
$ws.Range("D6").Value = 0.009174311926605505
$ws.Range("E6").Value = 0.003669724770642202
$ws.Range("F6").Value = 0.04954128440366973
$ws.Range("J6").Value = 0.2862385321100918
$ws.Range("O6").Value = 0.03119266055045872
$ws.Range("Q6").Value = 0.1541284403669725
$ws.Range("R6").Value = 0.04220183486238532
$ws.Range("S6").Value = 0.3486238532110092
$ws.Range("B7").Value = 0.1134453781512605
$ws.Range("D7").Value = 0.01260504201680672
$ws.Range("F7").Value = 0.05042016806722689
$ws.Range("J7").Value = 0.1890756302521008
$ws.Range("O7").Value = 0.02100840336134454
$ws.Range("Q7").Value = 0.1680672268907563
$ws.Range("R7").Value = 0.04831932773109244
$ws.Range("S7").Value = 0.3970588235294117
$ws.Range("B8").Value = 0.09982174688057041
$ws.Range("D8").Value = 0.01693404634581105
$ws.Range("E8").Value = 0.00089126559714795
$ws.Range("F8").Value = 0.06417112299465241
$ws.Range("J8").Value = 0.1167557932263815
$ws.Range("O8").Value = 0.02228163992869875
$ws.Range("Q8").Value = 0.1809269162210339
$ws.Range("R8").Value = 0.0659536541889483
$ws.Range("S8").Value = 0.4322638146167558
$ws.Range("B9").Value = 0.0968586387434555
$ws.Range("D9").Value = 0.005235602094240838
$ws.Range("E9").Value = 0.002617801047120419
$ws.Range("F9").Value = 0.08638743455497382
$ws.Range("J9").Value = 0.1230366492146597
$ws.Range("O9").Value = 0.01570680628272251
$ws.Range("Q9").Value = 0.143979057591623
$ws.Range("R9").Value = 0.08900523560209424
$ws.Range("S9").Value = 0.4371727748691099
$ws.Range("B10").Value = 0.1204819277108434
$ws.Range("D10").Value = 0.02371987951807229
$ws.Range("E10").Value = 0.0007530120481927711
$ws.Range("F10").Value = 0.07341867469879518
$ws.Range("J10").Value = 0.1065512048192771
$ws.Range("O10").Value = 0.0233433734939759
$ws.Range("Q10").Value = 0.2085843373493976
$ws.Range("R10").Value = 0.06626506024096386
$ws.Range("S10").Value = 0.376882530120482
$ws.Range("F11").Value = 0.001400560224089636
$ws.Range("G11").Value = 0.1666666666666667
$ws.Range("J11").Value = 0.09663865546218488
$ws.Range("K11").Value = 0.2030812324929972
$ws.Range("L11").Value = 0.5098039215686274
$ws.Range("S11").Value = 0.02240896358543417
$ws.Range("G12").Value = 0.7336448598130841
$ws.Range("J12").Value = 0.1682242990654206
$ws.Range("K12").Value = 0.01401869158878505
$ws.Range("L12").Value = 0.03271028037383177
$ws.Range("S12").Value = 0.0514018691588785
$ws.Range("F13").Value = 0.008928571428571428
$ws.Range("G13").Value = 0.5982142857142857
$ws.Range("J13").Value = 0.2857142857142857
$ws.Range("S13").Value = 0.1071428571428571
$ws.Range("G14").Value = 0.75
$ws.Range("J14").Value = 0.25
$ws.Range("F15").Value = 0.025390625
$ws.Range("H15").Value = 0.185546875
$ws.Range("I15").Value = 0.044921875
$ws.Range("J15").Value = 0.328125
$ws.Range("K15").Value = 0.072265625
$ws.Range("M15").Value = 0.01171875
$ws.Range("N15").Value = 0.00390625
$ws.Range("O15").Value = 0.060546875
$ws.Range("S15").Value = 0.267578125
$ws.Range("F16").Value = 0.02745995423340961
$ws.Range("H16").Value = 0.1624713958810069
$ws.Range("I16").Value = 0.07780320366132723
$ws.Range("J16").Value = 0.3821510297482837
$ws.Range("K16").Value = 0.1350114416475973
$ws.Range("M16").Value = 0.009153318077803204
$ws.Range("O16").Value = 0.07322654462242563
$ws.Range("S16").Value = 0.1327231121281464
$ws.Range("F17").Value = 0.02596053997923157
$ws.Range("H17").Value = 0.1962616822429906
$ws.Range("I17").Value = 0.0778816199376947
$ws.Range("J17").Value = 0.3935617860851506
$ws.Range("K17").Value = 0.102803738317757
$ws.Range("M17").Value = 0.02388369678089304
$ws.Range("N17").Value = 0.002076843198338525
$ws.Range("O17").Value = 0.06542056074766354
$ws.Range("S17").Value = 0.1121495327102804
$ws.Range("F18").Value = 0.03636363636363636
$ws.Range("H18").Value = 0.1757575757575758
$ws.Range("I18").Value = 0.09090909090909091
$ws.Range("J18").Value = 0.403030303030303
$ws.Range("K18").Value = 0.08484848484848485
$ws.Range("M18").Value = 0.02121212121212121
$ws.Range("O18").Value = 0.06060606060606061
$ws.Range("S18").Value = 0.1272727272727273
$ws.Range("F19").Value = 0.03132530120481928
$ws.Range("H19").Value = 0.2351118760757315
$ws.Range("I19").Value = 0.07091222030981068
$ws.Range("J19").Value = 0.3335628227194492
$ws.Range("K19").Value = 0.1070567986230637
$ws.Range("M19").Value = 0.023407917383821
$ws.Range("N19").Value = 0.001721170395869191
$ws.Range("O19").Value = 0.06506024096385542
$ws.Range("S19").Value = 0.13184165232358
